$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H39").Value = 209.3125
$ws.Range("I39").Value = 131.5
$ws.Range("J39").Value = 256
$ws.Range("K39").Value = 394.5
$ws.Range("L39").Value = 768
$ws.Range("M39").Value = -98.5
$ws.Range("N39").Value = -1360
$ws.Range("H40").Value = 1658.4615
$ws.Range("I40").Value = 1070
$ws.Range("J40").Value = 2600
$ws.Range("K40").Value = 1070
$ws.Range("L40").Value = 2600
$ws.Range("M40").Value = -895
$ws.Range("N40").Value = -2950
$ws.Range("H53").Value = 1491.5
$ws.Range("I53").Value = 275
$ws.Range("J53").Value = 1897
$ws.Range("K53").Value = 275
$ws.Range("L53").Value = 1897
$ws.Range("M53").Value = 362
$ws.Range("N53").Value = -3171
$ws.Range("H92").Value = 2250
$ws.Range("I92").Value = 2100
$ws.Range("K92").Value = 2100
$ws.Range("M92").Value = -852
$ws.Range("H97").Value = 1012
$ws.Range("J97").Value = 1012
$ws.Range("L97").Value = 3036
$ws.Range("N97").Value = -4028
$ws.Range("H101").Value = 1501.3334
$ws.Range("I101").Value = 1501.3334
$ws.Range("K101").Value = 4504.0002
$ws.Range("M101").Value = -2882.0002
$ws.Range("H105").Value = 37002.168
$ws.Range("J105").Value = 37002.168
$ws.Range("L105").Value = 37002.168
$ws.Range("N105").Value = -43990.168
$ws.Range("H116").Value = 6644.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 6644.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 6644.5
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -13528.5
$ws.Range("H133").Value = 52779.75
$ws.Range("J133").Value = 52779.75
$ws.Range("L133").Value = 52779.75
$ws.Range("N133").Value = -62899.75
$ws.Range("H136").Value = 42280
$ws.Range("J136").Value = 42280
$ws.Range("L136").Value = 42280
$ws.Range("N136").Value = -52480
$ws.Range("H137").Value = 92682.37
$ws.Range("I137").Value = 1643.5714
$ws.Range("J137").Value = 252000.25
$ws.Range("K137").Value = 4930.7142
$ws.Range("L137").Value = 756000.75
$ws.Range("M137").Value = -2380.7142
$ws.Range("N137").Value = -761100.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17312.016
$ws.Range("I32").Value = 17039.678
$ws.Range("K32").Value = 17039.678
$ws.Range("M32").Value = -16752.678
$ws.Range("H74").Value = 50003176
$ws.Range("I74").Value = 62503156
$ws.Range("J74").Value = 3250
$ws.Range("K74").Value = 62503156
$ws.Range("L74").Value = 3250
$ws.Range("M74").Value = -62502282
$ws.Range("N74").Value = -4998
$ws.Range("H77").Value = 50003176
$ws.Range("I77").Value = 62503156
$ws.Range("J77").Value = 3250
$ws.Range("K77").Value = 312515780
$ws.Range("L77").Value = 16250
$ws.Range("M77").Value = -312511412
$ws.Range("N77").Value = -24986
$ws.Range("H97").Value = 1285
$ws.Range("I97").Value = 1345.3846
$ws.Range("K97").Value = 1345.3846
$ws.Range("M97").Value = -849.3846000000001
$ws.Range("H102").Value = 1670.9
$ws.Range("I102").Value = 1451.6666
$ws.Range("J102").Value = 1999.75
$ws.Range("K102").Value = 1451.6666
$ws.Range("L102").Value = 1999.75
$ws.Range("M102").Value = 170.3334
$ws.Range("N102").Value = -5243.75
$ws.Range("H132").Value = 20666.26
$ws.Range("I132").Value = 1999
$ws.Range("J132").Value = 128003
$ws.Range("K132").Value = 5997
$ws.Range("L132").Value = 384009
$ws.Range("M132").Value = -3467
$ws.Range("N132").Value = -389069

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1453.7255
$ws.Range("I86").Value = 1377.1428
$ws.Range("J86").Value = 1811.1111
$ws.Range("K86").Value = 1377.1428
$ws.Range("L86").Value = 1811.1111
$ws.Range("M86").Value = -254.1428000000001
$ws.Range("N86").Value = -4057.1111
$ws.Range("H89").Value = 1453.7255
$ws.Range("I89").Value = 1377.1428
$ws.Range("J89").Value = 1811.1111
$ws.Range("K89").Value = 6885.714
$ws.Range("L89").Value = 9055.5555
$ws.Range("M89").Value = -1269.714
$ws.Range("N89").Value = -20287.5555
$ws.Range("H94").Value = 2020.341
$ws.Range("I94").Value = 1008.9091
$ws.Range("J94").Value = 5054.636
$ws.Range("K94").Value = 1008.9091
$ws.Range("L94").Value = 5054.636
$ws.Range("M94").Value = -557.9091
$ws.Range("N94").Value = -5956.636
$ws.Range("H99").Value = 1969.3125
$ws.Range("I99").Value = 1671.2858
$ws.Range("K99").Value = 1671.2858
$ws.Range("M99").Value = -173.2858000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 49999
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 49999
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -50471
$ws.Range("H30").Value = 49999
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 49999
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 49999
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -50181
$ws.Range("H31").Value = 3099.8
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3099.8
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3099.8
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3689.8
$ws.Range("H34").Value = 3099.8
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3099.8
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3099.8
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3503.8
$ws.Range("H58").Value = 18396.932
$ws.Range("I58").Value = 1231.8182
$ws.Range("J58").Value = 72344.42999999999
$ws.Range("K58").Value = 1231.8182
$ws.Range("L58").Value = 72344.42999999999
$ws.Range("M58").Value = -1028.8182
$ws.Range("N58").Value = -72750.42999999999
$ws.Range("H128").Value = 49999
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49999
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49999
$ws.Range("M128").ClearContents()
$ws.Range("N128").Value = -59959
$ws.Range("H136").Value = 18396.932
$ws.Range("I136").Value = 1231.8182
$ws.Range("J136").Value = 72344.42999999999
$ws.Range("K136").Value = 3695.4546
$ws.Range("L136").Value = 217033.29
$ws.Range("M136").Value = -1145.4546
$ws.Range("N136").Value = -222133.29

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1501.9
$ws.Range("I5").Value = 876.75
$ws.Range("J5").Value = 4002.5
$ws.Range("K5").Value = 2630.25
$ws.Range("L5").Value = 12007.5
$ws.Range("M5").Value = -2518.25
$ws.Range("N5").Value = -12231.5
$ws.Range("H36").Value = 3802.5
$ws.Range("I36").Value = 3802
$ws.Range("J36").Value = 3803
$ws.Range("K36").Value = 11406
$ws.Range("L36").Value = 11409
$ws.Range("M36").Value = -11237
$ws.Range("N36").Value = -11747
$ws.Range("H107").Value = 14473.857
$ws.Range("I107").Value = 33434.332
$ws.Range("J107").Value = 253.5
$ws.Range("K107").Value = 100302.996
$ws.Range("L107").Value = 760.5
$ws.Range("M107").Value = -98382.99600000001
$ws.Range("N107").Value = -4600.5
$ws.Range("H131").Value = 704.37
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 708.54083
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2125.62249
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12205.62249
$ws.Range("H132").Value = 1216.1333
$ws.Range("I132").Value = 1081
$ws.Range("J132").Value = 1370.5714
$ws.Range("K132").Value = 9729
$ws.Range("L132").Value = 12335.1426
$ws.Range("M132").Value = -7199
$ws.Range("N132").Value = -17395.1426
$ws.Range("H135").Value = 1501.9
$ws.Range("I135").Value = 876.75
$ws.Range("J135").Value = 4002.5
$ws.Range("K135").Value = 7890.75
$ws.Range("L135").Value = 36022.5
$ws.Range("M135").Value = -5355.75
$ws.Range("N135").Value = -41092.5
$ws.Range("H138").Value = 126497.71
$ws.Range("I138").Value = 1496.6666
$ws.Range("K138").Value = 4489.9998
$ws.Range("M138").Value = 650.0002000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3586.1794
$ws.Range("I126").Value = 2683.1538
$ws.Range("J126").Value = 5392.231
$ws.Range("K126").Value = 8049.4614
$ws.Range("L126").Value = 16176.693
$ws.Range("M126").Value = -5579.4614
$ws.Range("N126").Value = -21116.693

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 294.85184
$ws.Range("I16").Value = 294.29166
$ws.Range("J16").Value = 299.33334
$ws.Range("K16").Value = 294.29166
$ws.Range("L16").Value = 299.33334
$ws.Range("M16").Value = -124.29166
$ws.Range("N16").Value = -639.33334
$ws.Range("H136").Value = 39579.383
$ws.Range("I136").Value = 46593.816
$ws.Range("K136").Value = 139781.448
$ws.Range("M136").Value = -137231.448

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 40000
$ws.Range("J46").Value = 40000
$ws.Range("L46").Value = 40000
$ws.Range("N46").Value = -40462
$ws.Range("H107").Value = 1895079.4
$ws.Range("I107").Value = 811.9231
$ws.Range("J107").Value = 4133759
$ws.Range("K107").Value = 2435.7693
$ws.Range("L107").Value = 12401277
$ws.Range("M107").Value = -515.7692999999999
$ws.Range("N107").Value = -12405117
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -125070
$ws.Range("H135").Value = 52306
$ws.Range("J135").Value = 52306
$ws.Range("L135").Value = 52306
$ws.Range("N135").Value = -62446
